# Apply the "additional scraping" update to the PlayerPerformance workbook:
#  1. Insert a new "Player Info" sheet (as the first sheet) with the
#     player's ID, NAME, BATTING_HAND and BOWL_STYLE.
#  2. On the existing "ODI Batting" sheet, rename MATCH_CARD_LINK ->
#     MATCH_CODE and replace the full scorecard URL values with just the
#     numeric match code.
#  3. Do the same MATCH_CARD_LINK -> MATCH_CODE rename/value change on the
#     existing "ODI Bowling" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet -------------------------------------------
# Worksheets.Add() inserts the new sheet at the very front of the workbook,
# which is exactly where it needs to end up (before "ODI Batting").
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

# Keep the player ID as text (matches the inlineStr "5994" in the source).
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5994"
$playerInfo.Range("B2").Value = "Praveen Jayawickrama"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE -----------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4469", "4485", "4487", "4488", "4491")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $cell = $batting.Range("D" + (2 + $i))
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$i]
}

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE ------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4469", "4485", "4487", "4488", "4491")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $cell = $bowling.Range("B" + (2 + $i))
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$i]
}

Write-Host "Workbook updated: Player Info sheet added; MATCH_CARD_LINK -> MATCH_CODE on ODI Batting/Bowling."
